# Module 2: News Sentiment Analysis project
# Update predicted_sentiment (column L) and sentiment_score (column M) values
# on the single worksheet of the workbook to reflect corrected model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = -0.9
$ws.Range("M4").Value = 0.7
$ws.Range("M5").Value = 0.85
$ws.Range("M6").Value = -0.85
$ws.Range("M7").Value = -0.95
$ws.Range("M8").Value = -0.7
$ws.Range("M10").Value = -0.55
$ws.Range("L13").Value = "negative"
$ws.Range("M13").Value = -0.95
$ws.Range("L15").Value = "positive"
$ws.Range("M15").Value = 0.9
$ws.Range("L16").Value = "negative"
$ws.Range("M16").Value = -0.9
$ws.Range("L17").Value = "positive"
$ws.Range("M17").Value = 0.85
$ws.Range("L18").Value = "positive"
$ws.Range("M18").Value = 0.85
$ws.Range("L19").Value = "negative"
$ws.Range("M19").Value = -0.9
$ws.Range("L20").Value = "negative"
$ws.Range("M20").Value = -0.9
$ws.Range("L21").Value = "negative"
$ws.Range("M21").Value = -0.95
$ws.Range("L22").Value = "positive"
$ws.Range("M22").Value = 0.85
$ws.Range("L23").Value = "positive"
$ws.Range("M23").Value = 0.9
$ws.Range("L24").Value = "positive"
$ws.Range("M24").Value = 0.95
$ws.Range("M25").Value = -0.85
$ws.Range("L27").Value = "neutral"
$ws.Range("M27").Value = 0
$ws.Range("L28").Value = "neutral"
$ws.Range("M28").Value = 0
$ws.Range("L29").Value = "neutral"
$ws.Range("M29").Value = 0
$ws.Range("L30").Value = "neutral"
$ws.Range("M30").Value = 0
$ws.Range("L31").Value = "neutral"
$ws.Range("M31").Value = 0
$ws.Range("L32").Value = "neutral"
$ws.Range("M32").Value = 0
$ws.Range("L35").Value = "neutral"
$ws.Range("M35").Value = 0.15
$ws.Range("L36").Value = "positive"
$ws.Range("M36").Value = 0.85
$ws.Range("L38").Value = "negative"
$ws.Range("M38").Value = -0.5
$ws.Range("L39").Value = "negative"
$ws.Range("M39").Value = -0.95
$ws.Range("L40").Value = "positive"
$ws.Range("M40").Value = 0.9
$ws.Range("L41").Value = "positive"
$ws.Range("M41").Value = 0.75
$ws.Range("L42").Value = "positive"
$ws.Range("M42").Value = 0.8
$ws.Range("L43").Value = "positive"
$ws.Range("M43").Value = 0.9
$ws.Range("L44").Value = "negative"
$ws.Range("M44").Value = -0.85
$ws.Range("L45").Value = "positive"
$ws.Range("M45").Value = 0.95
$ws.Range("M47").Value = 0.98
$ws.Range("M48").Value = 1
$ws.Range("M52").Value = 0.9
$ws.Range("M55").Value = 0.4
$ws.Range("L56").Value = "neutral"
$ws.Range("M56").Value = 0
$ws.Range("L62").Value = "neutral"
$ws.Range("M62").Value = 0
$ws.Range("M63").Value = 0.9
$ws.Range("L64").Value = "neutral"
$ws.Range("M64").Value = 0
$ws.Range("L65").Value = "neutral"
$ws.Range("M65").Value = 0
$ws.Range("L66").Value = "neutral"
$ws.Range("M66").Value = 0
$ws.Range("M67").Value = 0.9
$ws.Range("M68").Value = -0.7
$ws.Range("M70").Value = 0.9
$ws.Range("M71").Value = -0.8
$ws.Range("L72").Value = "neutral"
$ws.Range("M72").Value = 0.1
$ws.Range("L73").Value = "positive"
$ws.Range("M73").Value = 0.95
$ws.Range("L74").Value = "positive"
$ws.Range("M74").Value = 0.98
$ws.Range("L75").Value = "negative"
$ws.Range("M75").Value = -0.8
$ws.Range("L76").Value = "negative"
$ws.Range("M76").Value = -0.95
$ws.Range("L84").Value = "positive"
$ws.Range("M84").Value = 0.7
$ws.Range("L85").Value = "negative"
$ws.Range("M85").Value = -0.95
$ws.Range("M86").Value = 0.35
$ws.Range("M88").Value = 0.9
$ws.Range("L90").Value = "positive"
$ws.Range("M90").Value = 0.5
$ws.Range("M91").Value = 0.85
$ws.Range("M92").Value = -0.35
$ws.Range("M93").Value = -0.85
$ws.Range("M94").Value = 0.8
$ws.Range("L97").Value = "positive"
$ws.Range("M97").Value = 0.6
